$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 3355.6310226535552
$ws.Range("C5").Value = 3648.0809990148869
$ws.Range("C6").Value = 3516.3657843565247
$ws.Range("C7").Value = 4404.0969892183703
$ws.Range("C8").Value = 4172.721419324309
$ws.Range("C9").Value = 5015.9554961040903
$ws.Range("C10").Value = 5439.0650793936511
$ws.Range("C11").Value = 5386.5079994594525
$ws.Range("C12").Value = 5392.4408464370508
$ws.Range("C13").Value = 5576.4507899787477
$ws.Range("C14").Value = 5485.2119506925355
$ws.Range("C15").Value = 5501.4965014230147
$ws.Range("C16").Value = 5579.2530405437765
$ws.Range("C17").Value = 5620.0178954415696
$ws.Range("C18").Value = 5558.1872775448492
$ws.Range("C19").Value = 5214.9476809986072
$ws.Range("C20").Value = 7886.7460001911013
$ws.Range("C21").Value = 4299.1560322753585
$ws.Range("C22").Value = 4177.7101415028701
$ws.Range("C23").Value = -5446.5476336899646
$ws.Range("C24").Value = 20336.627758370814
$ws.Range("C25").Value = 1020.9113378688274
$ws.Range("C26").Value = 877.39136048143951
$ws.Range("C27").Value = 890.5629436049453
$ws.Range("C28").Value = 855.89358197199385
$ws.Range("C29").Value = 818.56619953545169
$ws.Range("C30").Value = 1097.1827270047909
$ws.Range("C31").Value = 747.63784296569611
$ws.Range("C32").Value = 198.59772050578201
$ws.Range("C33").Value = -3189.1643075290062

$ws.Range("M10").Select()
